# Auto-generated Word COM-interop script to apply the diff
$d = $word.ActiveDocument

# Heading paragraph 22: open_character_stream (colored=True, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(22)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr><w:t>open_character_stream</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 27: open_token_stream (colored=True, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(27)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr><w:t>open_token_stream</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 32: get_token (colored=False, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(32)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/></w:pPr><w:r><w:t>get_token</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 44: is_token_end (colored=False, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(44)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/></w:pPr><w:r><w:t>is_token_end</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 55: token_type (colored=True, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(55)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr><w:t>token_type</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 66: print_token (colored=True, lastRenderedPageBreak=True)
$p = $d.Paragraphs.Item(66)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr><w:lastRenderedPageBreak/><w:t>print_token</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 71: is_comment (colored=True, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(71)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr><w:t>is_comment</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 76: is_keyword (colored=False, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(76)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/></w:pPr><w:r><w:t>is_keyword</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 81: is_char_constant (colored=True, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(81)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr><w:t>is_char_constant</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 86: is_num_constant (colored=False, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(86)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/></w:pPr><w:r><w:t>is_num_constant</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 92: is_str_constant (colored=False, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(92)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/></w:pPr><w:r><w:t>is_str_constant</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 98: is_identifier (colored=False, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(98)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/></w:pPr><w:r><w:t>is_identifier</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 104: print_spec_symbol (colored=True, lastRenderedPageBreak=True)
$p = $d.Paragraphs.Item(104)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr><w:lastRenderedPageBreak/><w:t>print_spec_symbol</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Heading paragraph 114: is_spec_symbol (colored=True, lastRenderedPageBreak=False)
$p = $d.Paragraphs.Item(114)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:keepLines w:val="0"/><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="C45911" w:themeColor="accent2" w:themeShade="BF"/></w:rPr><w:t>is_spec_symbol</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Paragraph 34: highlight 'Test Path 1: [Start, 1, 2, End]' under get_token
$p = $d.Paragraphs.Item(34)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:keepNext/></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Test </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Path</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> 1: [Start, 1, 2, End]</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Paragraph 124 (final, empty NoSpacing paragraph) -> 'main' heading + new test paths
$p = $d.Paragraphs.Item(124)
$p.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>main</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">Test Path 1: [main.Start, main.1, main.2, main.6, open_token_stream.Start, open_token_stream.1, open_token_stream.2, open_token_stream.4, open_token_stream.End, main.6, main.7, get_token.Start, get_token.1, get_token.2, get_token.End, main.7, main.8, main.9, print_token.Start, print_token.1, </w:t></w:r><w:r><w:t xml:space="preserve">token_type.Start, </w:t></w:r><w:r><w:t xml:space="preserve">, token_type.1, </w:t></w:r><w:r><w:t xml:space="preserve">is_keyword.Start, is_keyword.1, is_keyword.2, is_keyword.End, token_type.1, </w:t></w:r><w:r><w:t xml:space="preserve">token_type.2, token_type.End, </w:t></w:r><w:r><w:t xml:space="preserve">print_token.1, </w:t></w:r><w:r><w:t>print_token.2, print_token.3, print_token.4, print_token.5, print_token.6, print_token.7, print_token.8, print_token.9, print_token.10, print_token.11, print_token.12, print_token.End</w:t></w:r><w:r><w:t>, main.9, main.10,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>get_token.</w:t></w:r><w:r><w:t>Start, get_token.1, get_token.3, get_token.4, get_token.5, get_token.4, get_token.6, get_token.7, get_token.End</w:t></w:r><w:r><w:t>, main.10, main.End</w:t></w:r><w:r><w:t>]</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# The original trailing paragraph mark (now the last paragraph) becomes 'Test Path 2: '
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.InsertAfter("Test Path 2: ")

